# Apply the change described by the commit:
# [base commands] - [assertMatch(text,regex)]: NEW command to check for text
# value via regular expression.
#
# This edit touches the hidden '#system' reference sheet, which stores, per
# command-category column, the alphabetically sorted list of Nexial script
# commands (the column headers are the corresponding defined-name ranges):
#   1) A new "base" command  assertMatch(text,regex)   is inserted (column F)
#   2) A new "external" command openFile(filePath)     is inserted (column J)
#   3) The obsolete "tn.5250" category/column is removed altogether, which
#      means its name disappears from the "target" list (column A) and its
#      whole column (AA) is deleted, shifting everything after it left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Column F ("base"): insert "assertMatch(text,regex)" at F11, keeping
#    the alphabetical order, and push the remaining rows down by one
#    (old F11:F44 -> new F12:F45).
# ---------------------------------------------------------------------
$baseShift = $ws.Range("F11:F44").Value()
$ws.Range("F12:F45").Value = $baseShift
$ws.Range("F11").Value = "assertMatch(text,regex)"

# ---------------------------------------------------------------------
# 2) Column J ("external"): insert "openFile(filePath)" at J2, pushing the
#    remaining rows down by one (old J2:J6 -> new J3:J7).
# ---------------------------------------------------------------------
$externalShift = $ws.Range("J2:J6").Value()
$ws.Range("J3:J7").Value = $externalShift
$ws.Range("J2").Value = "openFile(filePath)"

# ---------------------------------------------------------------------
# 3) Column A ("target"): remove the "tn.5250" entry (row 27), pulling the
#    remaining category names up by one (old A28:A33 -> new A27:A32) and
#    clearing the now-unused last row.
# ---------------------------------------------------------------------
$targetShift = $ws.Range("A28:A33").Value()
$ws.Range("A27:A32").Value = $targetShift
$ws.Range("A33").ClearContents()

# ---------------------------------------------------------------------
# 4) Delete the whole "tn.5250" column (column AA, index 27); everything
#    to its right (web, webalert, webcookie, ws, ws.async, xml) shifts one
#    column to the left (AB->AA, AC->AB, AD->AC, AE->AD, AF->AE, AG->AF).
# ---------------------------------------------------------------------
$ws.Columns.Item(27).Delete()

# ---------------------------------------------------------------------
# 5) Update the defined names (workbook-level named ranges) so that they
#    keep pointing at the correct, resized/relocated ranges. Note that the
#    stale "tn.5250" name is intentionally left untouched/unremoved.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo        = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo    = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo      = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo         = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo    = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo   = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo          = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo    = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo         = "='#system'!`$AF`$2:`$AF`$27"

$wb.Save()
